# Auto-generated script to apply market-price / profit updates to the workbook
# across all 8 Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 25000
$ws.Range("N77").Value = -34360
$ws.Range("H80").Value = 914.2727
$ws.Range("J80").Value = 1140.6
$ws.Range("L80").Value = 3421.8
$ws.Range("N80").Value = -5417.799999999999
$ws.Range("H83").Value = 914.2727
$ws.Range("J83").Value = 1140.6
$ws.Range("L83").Value = 10265.4
$ws.Range("N83").Value = -20249.4
$ws.Range("H92").Value = 3750.75
$ws.Range("I92").Value = 1668
$ws.Range("K92").Value = 1668
$ws.Range("M92").Value = -420
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492
$ws.Range("H137").Value = 10292
$ws.Range("I137").Value = 7702.857
$ws.Range("K137").Value = 23108.571
$ws.Range("M137").Value = -20558.571
$ws.Range("H138").Value = 1000
$ws.Range("I138").Value = 1000
$ws.Range("K138").Value = 3000
$ws.Range("M138").Value = 2140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14318.75
$ws.Range("I28").Value = 14318.75
$ws.Range("K28").Value = 14318.75
$ws.Range("M28").Value = -14126.75
$ws.Range("H88").Value = 1515.1
$ws.Range("I88").Value = 1492
$ws.Range("J88").Value = 1549.75
$ws.Range("K88").Value = 1492
$ws.Range("L88").Value = 1549.75
$ws.Range("M88").Value = -1086
$ws.Range("N88").Value = -2361.75
$ws.Range("H91").Value = 1515.1
$ws.Range("I91").Value = 1492
$ws.Range("J91").Value = 1549.75
$ws.Range("K91").Value = 1492
$ws.Range("L91").Value = 1549.75
$ws.Range("M91").Value = -88
$ws.Range("N91").Value = -4357.75
$ws.Range("H97").Value = 1801
$ws.Range("I97").Value = 1801
$ws.Range("K97").Value = 1801
$ws.Range("M97").Value = -1305
$ws.Range("H99").Value = 14318.75
$ws.Range("I99").Value = 14318.75
$ws.Range("K99").Value = 14318.75
$ws.Range("M99").Value = -11323.75
$ws.Range("H102").Value = 1658
$ws.Range("I102").Value = 1658
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1658
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -36
$ws.Range("N102").Value = ""
$ws.Range("H122").Value = 3666.6667
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H132").Value = 6874
$ws.Range("I132").Value = 1148
$ws.Range("J132").Value = 12600
$ws.Range("K132").Value = 3444
$ws.Range("L132").Value = 37800
$ws.Range("M132").Value = -914
$ws.Range("N132").Value = -42860
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 249999
$ws.Range("J42").Value = 249999
$ws.Range("L42").Value = 249999
$ws.Range("N42").Value = -250655
$ws.Range("H54").Value = 2027.6666
$ws.Range("I54").Value = 2027.6666
$ws.Range("K54").Value = 2027.6666
$ws.Range("M54").Value = -1543.6666
$ws.Range("H94").Value = 2166.6667
$ws.Range("I94").Value = 2166.6667
$ws.Range("K94").Value = 2166.6667
$ws.Range("M94").Value = -1715.6667
$ws.Range("H97").Value = 6114
$ws.Range("I97").Value = 6114
$ws.Range("K97").Value = 6114
$ws.Range("M97").Value = -5123
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -1002
$ws.Range("N99").Value = -5496
$ws.Range("H105").Value = 2625.75
$ws.Range("I105").Value = 1156.4286
$ws.Range("K105").Value = 1156.4286
$ws.Range("M105").Value = 590.5714
$ws.Range("H107").Value = 527.375
$ws.Range("I107").Value = 501.2857
$ws.Range("J107").Value = 710
$ws.Range("K107").Value = 501.2857
$ws.Range("L107").Value = 710
$ws.Range("M107").Value = 1418.7143
$ws.Range("N107").Value = -4550

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 9566.333000000001
$ws.Range("J95").Value = 9566.333000000001
$ws.Range("L95").Value = 9566.333000000001
$ws.Range("N95").Value = -15058.333
$ws.Range("H96").Value = 7871.6665
$ws.Range("J96").Value = 7871.6665
$ws.Range("L96").Value = 7871.6665
$ws.Range("N96").Value = -13363.6665
$ws.Range("H105").Value = 1992.5
$ws.Range("I105").Value = 1992.5
$ws.Range("K105").Value = 1992.5
$ws.Range("M105").Value = -245.5
$ws.Range("H134").Value = 6002.3335
$ws.Range("I134").Value = 1601.4
$ws.Range("K134").Value = 4804.200000000001
$ws.Range("M134").Value = -2269.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 400
$ws.Range("J3").Value = 400
$ws.Range("L3").Value = 1200
$ws.Range("N3").Value = -1424
$ws.Range("H11").Value = 612
$ws.Range("J11").Value = 999
$ws.Range("L11").Value = 2997
$ws.Range("N11").Value = -3277
$ws.Range("H21").Value = 300
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 900
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -1246
$ws.Range("H33").Value = 951
$ws.Range("I33").Value = 872.5714
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 5235.428400000001
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = -4952.428400000001
$ws.Range("N33").Value = -9566
$ws.Range("H34").Value = 1868
$ws.Range("J34").Value = 2996.6667
$ws.Range("L34").Value = 8990.000100000001
$ws.Range("N34").Value = -9158.000100000001
$ws.Range("H81").Value = 513
$ws.Range("I81").Value = 513
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1539
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -416
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 513
$ws.Range("I84").Value = 513
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4617
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 999
$ws.Range("N84").Value = ""
$ws.Range("H98").Value = 617
$ws.Range("I98").Value = 852
$ws.Range("J98").Value = 147
$ws.Range("K98").Value = 2556
$ws.Range("L98").Value = 441
$ws.Range("M98").Value = -1058
$ws.Range("N98").Value = -3437
$ws.Range("H108").Value = 165.75
$ws.Range("I108").Value = 165.75
$ws.Range("K108").Value = 497.25
$ws.Range("M108").Value = 2382.75
$ws.Range("H117").Value = 2348.5
$ws.Range("I117").Value = 1964.6666
$ws.Range("J117").Value = 3500
$ws.Range("K117").Value = 5893.9998
$ws.Range("L117").Value = 10500
$ws.Range("M117").Value = -2451.9998
$ws.Range("N117").Value = -17384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2602.5
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 5005
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 5005
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -5229
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 1000
$ws.Range("M70").Value = -730
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 1000
$ws.Range("M73").Value = -64
$ws.Range("H80").Value = 2199.6
$ws.Range("I80").Value = 1666.3334
$ws.Range("J80").Value = 2999.5
$ws.Range("K80").Value = 1666.3334
$ws.Range("L80").Value = 2999.5
$ws.Range("M80").Value = -668.3334
$ws.Range("N80").Value = -4995.5
$ws.Range("H83").Value = 2199.6
$ws.Range("I83").Value = 1666.3334
$ws.Range("J83").Value = 2999.5
$ws.Range("K83").Value = 8331.666999999999
$ws.Range("L83").Value = 14997.5
$ws.Range("M83").Value = -3339.666999999999
$ws.Range("N83").Value = -24981.5
$ws.Range("H97").Value = 514.8333
$ws.Range("I97").Value = 424.75
$ws.Range("K97").Value = 424.75
$ws.Range("M97").Value = 71.25
$ws.Range("H114").Value = 47500
$ws.Range("J114").Value = 47500
$ws.Range("L114").Value = 47500
$ws.Range("N114").Value = -56178
$ws.Range("H118").Value = 20899
$ws.Range("J118").Value = 20899
$ws.Range("L118").Value = 20899
$ws.Range("N118").Value = -24213
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 6044.923
$ws.Range("I132").Value = 3310.2222
$ws.Range("K132").Value = 9930.6666
$ws.Range("M132").Value = -7400.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H100").Value = 2584.5715
$ws.Range("I100").Value = 3162.6667
$ws.Range("J100").Value = 2151
$ws.Range("K100").Value = 3162.6667
$ws.Range("L100").Value = 2151
$ws.Range("M100").Value = -2621.6667
$ws.Range("N100").Value = -3233
$ws.Range("H136").Value = 25000
$ws.Range("J136").Value = 25000
$ws.Range("L136").Value = 75000
$ws.Range("N136").Value = -80100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 35285.5
$ws.Range("J97").Value = 35285.5
$ws.Range("L97").Value = 35285.5
$ws.Range("N97").Value = -37267.5
$ws.Range("H136").Value = 12987.667
$ws.Range("I136").Value = 8472.5
$ws.Range("K136").Value = 25417.5
$ws.Range("M136").Value = -22867.5
